$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# Row 21
$ws.Cells.Item(21, 8).Value = 10950
$ws.Cells.Item(21, 9).Value = 12000
$ws.Cells.Item(21, 10).Value = 9900
$ws.Cells.Item(21, 11).Value = 12000
$ws.Cells.Item(21, 12).Value = 9900
$ws.Cells.Item(21, 13).Value = -11532
$ws.Cells.Item(21, 14).Value = -10836

# Row 23
$ws.Cells.Item(23, 8).Value = 10950
$ws.Cells.Item(23, 9).Value = 12000
$ws.Cells.Item(23, 10).Value = 9900
$ws.Cells.Item(23, 11).Value = 12000
$ws.Cells.Item(23, 12).Value = 9900
$ws.Cells.Item(23, 13).Value = -11766
$ws.Cells.Item(23, 14).Value = -10368

# Row 70
$ws.Cells.Item(70, 8).Value = 857.3333
$ws.Cells.Item(70, 9).Value = 671.3333
$ws.Cells.Item(70, 10).Value = 1043.3334
$ws.Cells.Item(70, 11).Value = 2013.9999
$ws.Cells.Item(70, 12).Value = 3130.0002
$ws.Cells.Item(70, 13).Value = -1743.9999
$ws.Cells.Item(70, 14).Value = -3670.0002

# Row 73
$ws.Cells.Item(73, 8).Value = 857.3333
$ws.Cells.Item(73, 9).Value = 671.3333
$ws.Cells.Item(73, 10).Value = 1043.3334
$ws.Cells.Item(73, 11).Value = 2013.9999
$ws.Cells.Item(73, 12).Value = 3130.0002
$ws.Cells.Item(73, 13).Value = -1077.9999
$ws.Cells.Item(73, 14).Value = -5002.0002

# Row 76
$ws.Cells.Item(76, 8).Value = 85745.44500000001
$ws.Cells.Item(76, 9).Value = 113911.15
$ws.Cells.Item(76, 10).Value = 5272
$ws.Cells.Item(76, 11).Value = 113911.15
$ws.Cells.Item(76, 12).Value = 5272
$ws.Cells.Item(76, 13).Value = -113596.15
$ws.Cells.Item(76, 14).Value = -5902

# Row 79
$ws.Cells.Item(79, 8).Value = 85745.44500000001
$ws.Cells.Item(79, 9).Value = 113911.15
$ws.Cells.Item(79, 10).Value = 5272
$ws.Cells.Item(79, 11).Value = 113911.15
$ws.Cells.Item(79, 12).Value = 5272
$ws.Cells.Item(79, 13).Value = -112819.15
$ws.Cells.Item(79, 14).Value = -7456

# Row 82
$ws.Cells.Item(82, 8).Value = 1704.8
$ws.Cells.Item(82, 9).Value = 507.25
$ws.Cells.Item(82, 10).Value = 6495
$ws.Cells.Item(82, 11).Value = 1521.75
$ws.Cells.Item(82, 12).Value = 19485
$ws.Cells.Item(82, 13).Value = -1115.75
$ws.Cells.Item(82, 14).Value = -20297

# Row 85
$ws.Cells.Item(85, 8).Value = 1704.8
$ws.Cells.Item(85, 9).Value = 507.25
$ws.Cells.Item(85, 10).Value = 6495
$ws.Cells.Item(85, 11).Value = 1521.75
$ws.Cells.Item(85, 12).Value = 19485
$ws.Cells.Item(85, 13).Value = -117.75
$ws.Cells.Item(85, 14).Value = -22293

# Row 88
$ws.Cells.Item(88, 8).Value = 2284.394
$ws.Cells.Item(88, 9).Value = 5054.8887
$ws.Cells.Item(88, 10).Value = 1245.4584
$ws.Cells.Item(88, 11).Value = 5054.8887
$ws.Cells.Item(88, 12).Value = 1245.4584
$ws.Cells.Item(88, 13).Value = -4648.8887
$ws.Cells.Item(88, 14).Value = -2057.4584

# Row 91
$ws.Cells.Item(91, 8).Value = 2284.394
$ws.Cells.Item(91, 9).Value = 5054.8887
$ws.Cells.Item(91, 10).Value = 1245.4584
$ws.Cells.Item(91, 11).Value = 5054.8887
$ws.Cells.Item(91, 12).Value = 1245.4584
$ws.Cells.Item(91, 13).Value = -3650.8887
$ws.Cells.Item(91, 14).Value = -4053.4584

# Row 97
$ws.Cells.Item(97, 8).Value = 285715600
$ws.Cells.Item(97, 9).Value = 200001200
$ws.Cells.Item(97, 10).Value = 500001500
$ws.Cells.Item(97, 11).Value = 600003600
$ws.Cells.Item(97, 12).Value = 1500004500
$ws.Cells.Item(97, 13).Value = -600003104
$ws.Cells.Item(97, 14).Value = -1500005492

# Row 100
$ws.Cells.Item(100, 8).Value = 2086.1482
$ws.Cells.Item(100, 9).Value = 1147.6923
$ws.Cells.Item(100, 10).Value = 2957.5715
$ws.Cells.Item(100, 11).Value = 1147.6923
$ws.Cells.Item(100, 12).Value = 2957.5715
$ws.Cells.Item(100, 13).Value = -606.6922999999999
$ws.Cells.Item(100, 14).Value = -4039.5715

# Row 103
$ws.Cells.Item(103, 8).Value = 465.27274
$ws.Cells.Item(103, 9).Value = 491
$ws.Cells.Item(103, 10).Value = 349.5
$ws.Cells.Item(103, 11).Value = 1473
$ws.Cells.Item(103, 12).Value = 1048.5
$ws.Cells.Item(103, 13).Value = -887
$ws.Cells.Item(103, 14).Value = -2220.5

# Row 106
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).ClearContents()
$ws.Cells.Item(106, 14).ClearContents()

# Row 109
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 14).ClearContents()

# Row 112
$ws.Cells.Item(112, 8).Value = 1116906.8
$ws.Cells.Item(112, 9).Value = 549.75
$ws.Cells.Item(112, 10).Value = 1563449.5
$ws.Cells.Item(112, 11).Value = 1649.25
$ws.Cells.Item(112, 12).Value = 4690348.5
$ws.Cells.Item(112, 13).Value = -541.25
$ws.Cells.Item(112, 14).Value = -4692564.5

# Row 115
$ws.Cells.Item(115, 8).Value = 381.15384
$ws.Cells.Item(115, 9).Value = 381.15384
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 1143.46152
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 423.5384799999999
$ws.Cells.Item(115, 14).ClearContents()

# Row 118
$ws.Cells.Item(118, 8).Value = 4441.222
$ws.Cells.Item(118, 9).Value = 598.2222
$ws.Cells.Item(118, 10).Value = 8284.223
$ws.Cells.Item(118, 11).Value = 1794.6666
$ws.Cells.Item(118, 12).Value = 24852.669
$ws.Cells.Item(118, 13).Value = -137.6666
$ws.Cells.Item(118, 14).Value = -28166.669

# Row 131
$ws.Cells.Item(131, 8).Value = 1923.75
$ws.Cells.Item(131, 9).Value = 565
$ws.Cells.Item(131, 10).Value = 6000
$ws.Cells.Item(131, 11).Value = 1695
$ws.Cells.Item(131, 12).Value = 18000
$ws.Cells.Item(131, 13).Value = 3345
$ws.Cells.Item(131, 14).Value = -28080

# Row 138
$ws.Cells.Item(138, 8).Value = 5425.36
$ws.Cells.Item(138, 9).Value = 1484.5312
$ws.Cells.Item(138, 10).Value = 7279.8677
$ws.Cells.Item(138, 11).Value = 4453.5936
$ws.Cells.Item(138, 12).Value = 21839.6031
$ws.Cells.Item(138, 13).Value = 686.4063999999998
$ws.Cells.Item(138, 14).Value = -32119.6031

$ws = $wb.Sheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 3516.5
$ws.Cells.Item(61, 9).Value = 3219.8
$ws.Cells.Item(61, 10).Value = 5000
$ws.Cells.Item(61, 11).Value = 3219.8
$ws.Cells.Item(61, 12).Value = 5000
$ws.Cells.Item(61, 13).Value = -3007.8
$ws.Cells.Item(61, 14).Value = -5424

# Row 75
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()

# Row 78
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()

# Row 136
$ws.Cells.Item(136, 8).Value = 3516.5
$ws.Cells.Item(136, 9).Value = 3219.8
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 9659.400000000001
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -7109.400000000001
$ws.Cells.Item(136, 14).Value = -20100

$ws = $wb.Sheets.Item("CUL")
# Row 92
$ws.Cells.Item(92, 8).Value = 520.8889
$ws.Cells.Item(92, 9).Value = 377.8
$ws.Cells.Item(92, 10).Value = 699.75
$ws.Cells.Item(92, 11).Value = 1133.4
$ws.Cells.Item(92, 12).Value = 2099.25
$ws.Cells.Item(92, 13).Value = 114.5999999999999
$ws.Cells.Item(92, 14).Value = -4595.25
